# Applies "major fixes to fields" commit to LCA_infrastructure.xlsx
# Touches: DHW, HEATING, COOLING, ELECTRICITY sheets.

$wb = $excel.ActiveWorkbook

$dhw = $wb.Worksheets.Item("DHW")
$heating = $wb.Worksheets.Item("HEATING")
$cooling = $wb.Worksheets.Item("COOLING")
$electricity = $wb.Worksheets.Item("ELECTRICITY")

# ---------------------------------------------------------------------------
# 1. Shared-string value fixes (SOURCE column): "FUEL" -> "BOILER",
#    "RENEWABLE" -> "SC" on both DHW and HEATING sheets.
# ---------------------------------------------------------------------------

foreach ($ws in @($dhw, $heating)) {
    foreach ($r in @("C3", "C4", "C5", "C7")) {
        $ws.Range($r).Value = "BOILER"
    }
    foreach ($r in @("C9", "C16", "C17")) {
        $ws.Range($r).Value = "SC"
    }
}

# ---------------------------------------------------------------------------
# 2. Selection / active-cell bookkeeping (cosmetic, matches authored file).
# ---------------------------------------------------------------------------

$dhw.Range("E1:E1048576").Select()
$heating.Range("E1:E1048576").Select()
$cooling.Range("C2").Select()

# ---------------------------------------------------------------------------
# 3. HEATING formatting fixes - align styles of E1, E7, E17:E21 with the
#    equivalent already-fixed cells on the DHW sheet.
# ---------------------------------------------------------------------------

$dhw.Range("D1").Copy() | Out-Null
$heating.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$dhw.Range("E7").Copy() | Out-Null
$heating.Range("E7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

foreach ($r in @("E17", "E18", "E19", "E20", "E21")) {
    $dhw.Range($r).Copy() | Out-Null
    $heating.Range($r).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. ELECTRICITY!C8 now derives from HEATING!D5 instead of HEATING!E5.
# ---------------------------------------------------------------------------

$electricity.Range("C8").Formula = "=HEATING!D5/0.4"

# COOLING!D6 keeps its formula (ELECTRICITY!C8/4); its cached value will be
# refreshed automatically once Excel recalculates the workbook.
$excel.CalculateFullRebuild()

# Re-select the final active sheet/cell to mirror the authored workbook
# (HEATING tab was the active tab before and after the edit).
$heating.Activate()
